$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = 'Última actualización: 07:40:11'
$ws1.Cells.Item(3,1).Value = 'Total filas: 60'

$ws1.Cells.Item(37,1).Value = '07:40:11'
$ws1.Cells.Item(37,2).Value = '07:44'
$ws1.Cells.Item(37,3).Value = '215A_EL PATO'
$ws1.Cells.Item(37,4).Value = 4
$ws1.Cells.Item(37,5).Value = 'LP1912'

$ws1.Cells.Item(38,1).Value = '06:38:54'
$ws1.Cells.Item(38,2).Value = '07:54'
$ws1.Cells.Item(38,3).Value = '14_ABASTO'
$ws1.Cells.Item(38,4).Value = 76
$ws1.Cells.Item(38,5).Value = 'LP1912'

$ws1.Cells.Item(39,1).Value = '07:40:11'
$ws1.Cells.Item(39,2).Value = '07:55'
$ws1.Cells.Item(39,3).Value = '14_ABASTO'
$ws1.Cells.Item(39,4).Value = 15
$ws1.Cells.Item(39,5).Value = 'LP1912'

$ws1.Cells.Item(40,1).Value = '06:19:59'
$ws1.Cells.Item(40,2).Value = '07:59'
$ws1.Cells.Item(40,3).Value = '17_ROMERO'
$ws1.Cells.Item(40,4).Value = 100
$ws1.Cells.Item(40,5).Value = 'LP1912'

$ws1.Cells.Item(41,1).Value = '06:38:54'
$ws1.Cells.Item(41,2).Value = '08:00'
$ws1.Cells.Item(41,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(41,4).Value = 82
$ws1.Cells.Item(41,5).Value = 'LP1912'

$ws1.Cells.Item(42,1).Value = '07:40:11'
$ws1.Cells.Item(42,2).Value = '08:00'
$ws1.Cells.Item(42,3).Value = '17_ROMERO'
$ws1.Cells.Item(42,4).Value = 20
$ws1.Cells.Item(42,5).Value = 'LP1912'

$ws1.Cells.Item(43,1).Value = '07:40:11'
$ws1.Cells.Item(43,2).Value = '08:01'
$ws1.Cells.Item(43,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(43,4).Value = 21
$ws1.Cells.Item(43,5).Value = 'LP1912'

$ws1.Cells.Item(44,1).Value = '07:40:11'
$ws1.Cells.Item(44,2).Value = '08:06'
$ws1.Cells.Item(44,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(44,4).Value = 26
$ws1.Cells.Item(44,5).Value = 'LP1912'

$ws1.Cells.Item(45,1).Value = '07:40:11'
$ws1.Cells.Item(45,2).Value = '08:11'
$ws1.Cells.Item(45,3).Value = '10_OLMOS'
$ws1.Cells.Item(45,4).Value = 31
$ws1.Cells.Item(45,5).Value = 'LP1912'

$ws1.Cells.Item(46,1).Value = '06:19:59'
$ws1.Cells.Item(46,2).Value = '08:12'
$ws1.Cells.Item(46,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(46,4).Value = 113
$ws1.Cells.Item(46,5).Value = 'LP1912'

$ws1.Cells.Item(47,1).Value = '07:40:11'
$ws1.Cells.Item(47,2).Value = '08:13'
$ws1.Cells.Item(47,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(47,4).Value = 33
$ws1.Cells.Item(47,5).Value = 'LP1912'

$ws1.Cells.Item(48,1).Value = '06:38:54'
$ws1.Cells.Item(48,2).Value = '08:28'
$ws1.Cells.Item(48,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(48,4).Value = 110
$ws1.Cells.Item(48,5).Value = 'LP1912'

$ws1.Cells.Item(49,1).Value = '07:40:11'
$ws1.Cells.Item(49,2).Value = '08:29'
$ws1.Cells.Item(49,3).Value = '15_ABASTO'
$ws1.Cells.Item(49,4).Value = 49
$ws1.Cells.Item(49,5).Value = 'LP1912'

$ws1.Cells.Item(50,1).Value = '07:40:11'
$ws1.Cells.Item(50,2).Value = '08:29'
$ws1.Cells.Item(50,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(50,4).Value = 49
$ws1.Cells.Item(50,5).Value = 'LP1912'

$ws1.Cells.Item(51,1).Value = '07:40:11'
$ws1.Cells.Item(51,2).Value = '08:41'
$ws1.Cells.Item(51,3).Value = '10_OLMOS'
$ws1.Cells.Item(51,4).Value = 61
$ws1.Cells.Item(51,5).Value = 'LP1912'

$ws1.Cells.Item(52,1).Value = '07:40:11'
$ws1.Cells.Item(52,2).Value = '08:41'
$ws1.Cells.Item(52,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(52,4).Value = 61
$ws1.Cells.Item(52,5).Value = 'LP1912'

$ws1.Cells.Item(53,1).Value = '07:15:48'
$ws1.Cells.Item(53,2).Value = '08:43'
$ws1.Cells.Item(53,3).Value = '215C_EL PATO'
$ws1.Cells.Item(53,4).Value = 88
$ws1.Cells.Item(53,5).Value = 'LP1912'

$ws1.Cells.Item(54,1).Value = '07:40:11'
$ws1.Cells.Item(54,2).Value = '08:44'
$ws1.Cells.Item(54,3).Value = '215C_EL PATO'
$ws1.Cells.Item(54,4).Value = 64
$ws1.Cells.Item(54,5).Value = 'LP1912'

$ws1.Cells.Item(55,1).Value = '07:40:11'
$ws1.Cells.Item(55,2).Value = '08:46'
$ws1.Cells.Item(55,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(55,4).Value = 66
$ws1.Cells.Item(55,5).Value = 'LP1912'

$ws1.Cells.Item(56,1).Value = '07:15:48'
$ws1.Cells.Item(56,2).Value = '08:51'
$ws1.Cells.Item(56,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(56,4).Value = 96
$ws1.Cells.Item(56,5).Value = 'LP1912'

$ws1.Cells.Item(57,1).Value = '06:56:24'
$ws1.Cells.Item(57,2).Value = '08:52'
$ws1.Cells.Item(57,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(57,4).Value = 116
$ws1.Cells.Item(57,5).Value = 'LP1912'

$ws1.Cells.Item(58,1).Value = '07:15:48'
$ws1.Cells.Item(58,2).Value = '08:53'
$ws1.Cells.Item(58,3).Value = '215B_EL PATO'
$ws1.Cells.Item(58,4).Value = 98
$ws1.Cells.Item(58,5).Value = 'LP1912'

$ws1.Cells.Item(59,1).Value = '07:40:11'
$ws1.Cells.Item(59,2).Value = '08:54'
$ws1.Cells.Item(59,3).Value = '215B_EL PATO'
$ws1.Cells.Item(59,4).Value = 74
$ws1.Cells.Item(59,5).Value = 'LP1912'

$ws1.Cells.Item(60,1).Value = '07:15:48'
$ws1.Cells.Item(60,2).Value = '08:57'
$ws1.Cells.Item(60,3).Value = '215A_EL PATO'
$ws1.Cells.Item(60,4).Value = 102
$ws1.Cells.Item(60,5).Value = 'LP1912'

$ws1.Cells.Item(61,1).Value = '07:40:11'
$ws1.Cells.Item(61,2).Value = '08:58'
$ws1.Cells.Item(61,3).Value = '215A_EL PATO'
$ws1.Cells.Item(61,4).Value = 78
$ws1.Cells.Item(61,5).Value = 'LP1912'

$ws1.Cells.Item(62,1).Value = '07:40:11'
$ws1.Cells.Item(62,2).Value = '09:14'
$ws1.Cells.Item(62,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(62,4).Value = 94
$ws1.Cells.Item(62,5).Value = 'LP1912'

$ws1.Cells.Item(63,1).Value = '07:40:11'
$ws1.Cells.Item(63,2).Value = '09:18'
$ws1.Cells.Item(63,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(63,4).Value = 98
$ws1.Cells.Item(63,5).Value = 'LP1912'

$ws1.Cells.Item(64,1).Value = '07:40:11'
$ws1.Cells.Item(64,2).Value = '09:18'
$ws1.Cells.Item(64,3).Value = '14_ABASTO'
$ws1.Cells.Item(64,4).Value = 98
$ws1.Cells.Item(64,5).Value = 'LP1912'

$ws1.Cells.Item(65,1).Value = '07:40:11'
$ws1.Cells.Item(65,2).Value = '09:31'
$ws1.Cells.Item(65,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(65,4).Value = 111
$ws1.Cells.Item(65,5).Value = 'LP1912'

# ---------- Sheet 2: LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = 'Última actualización: 07:40:11'
$ws2.Cells.Item(3,1).Value = 'Total filas: 13'

$ws2.Cells.Item(12,1).Value = '07:40:11'
$ws2.Cells.Item(12,2).Value = '07:44'
$ws2.Cells.Item(12,3).Value = '215A_EL PATO'
$ws2.Cells.Item(12,4).Value = 4
$ws2.Cells.Item(12,5).Value = 'LP1912'

$ws2.Cells.Item(13,1).Value = '07:15:48'
$ws2.Cells.Item(13,2).Value = '08:43'
$ws2.Cells.Item(13,3).Value = '215C_EL PATO'
$ws2.Cells.Item(13,4).Value = 88
$ws2.Cells.Item(13,5).Value = 'LP1912'

$ws2.Cells.Item(14,1).Value = '07:40:11'
$ws2.Cells.Item(14,2).Value = '08:44'
$ws2.Cells.Item(14,3).Value = '215C_EL PATO'
$ws2.Cells.Item(14,4).Value = 64
$ws2.Cells.Item(14,5).Value = 'LP1912'

$ws2.Cells.Item(15,1).Value = '07:15:48'
$ws2.Cells.Item(15,2).Value = '08:53'
$ws2.Cells.Item(15,3).Value = '215B_EL PATO'
$ws2.Cells.Item(15,4).Value = 98
$ws2.Cells.Item(15,5).Value = 'LP1912'

$ws2.Cells.Item(16,1).Value = '07:40:11'
$ws2.Cells.Item(16,2).Value = '08:54'
$ws2.Cells.Item(16,3).Value = '215B_EL PATO'
$ws2.Cells.Item(16,4).Value = 74
$ws2.Cells.Item(16,5).Value = 'LP1912'

$ws2.Cells.Item(17,1).Value = '07:15:48'
$ws2.Cells.Item(17,2).Value = '08:57'
$ws2.Cells.Item(17,3).Value = '215A_EL PATO'
$ws2.Cells.Item(17,4).Value = 102
$ws2.Cells.Item(17,5).Value = 'LP1912'

$ws2.Cells.Item(18,1).Value = '07:40:11'
$ws2.Cells.Item(18,2).Value = '08:58'
$ws2.Cells.Item(18,3).Value = '215A_EL PATO'
$ws2.Cells.Item(18,4).Value = 78
$ws2.Cells.Item(18,5).Value = 'LP1912'

# ---------- Sheet 3: 6203-6173 ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = 'Última actualización: 07:40:11'

$ws3.Cells.Item(6,1).Value = '07:15:48'
$ws3.Cells.Item(6,2).Value = '07:42'
$ws3.Cells.Item(6,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(6,4).Value = 27
$ws3.Cells.Item(6,5).Value = 'L6173'

$ws3.Cells.Item(7,1).Value = '07:40:11'
$ws3.Cells.Item(7,2).Value = '07:43'
$ws3.Cells.Item(7,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(7,4).Value = 3
$ws3.Cells.Item(7,5).Value = 'L6173'

$ws3.Cells.Item(8,1).Value = '07:15:48'
$ws3.Cells.Item(8,2).Value = '08:35'
$ws3.Cells.Item(8,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(8,4).Value = 80
$ws3.Cells.Item(8,5).Value = 'L6173'

$ws3.Cells.Item(9,1).Value = '07:40:11'
$ws3.Cells.Item(9,2).Value = '08:36'
$ws3.Cells.Item(9,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(9,4).Value = 56
$ws3.Cells.Item(9,5).Value = 'L6173'

$ws3.Cells.Item(10,1).Value = '07:15:48'
$ws3.Cells.Item(10,2).Value = '08:50'
$ws3.Cells.Item(10,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(10,4).Value = 95
$ws3.Cells.Item(10,5).Value = 'L6203'

$ws3.Cells.Item(11,1).Value = '07:40:11'
$ws3.Cells.Item(11,2).Value = '08:51'
$ws3.Cells.Item(11,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(11,4).Value = 71
$ws3.Cells.Item(11,5).Value = 'L6203'

Write-Output "done"